$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1) Rows 19 / 20 and 84 / 85 had their match-detail columns (F:V)
#    swapped between the row pair (index/date columns A:E stay put).
# -----------------------------------------------------------------
function Swap-RowRange($row1, $row2, $colStart, $colEnd) {
    $vals1 = @{}
    $vals2 = @{}
    for ($c = $colStart; $c -le $colEnd; $c++) {
        $vals1[$c] = $ws.Cells.Item($row1, $c).Value2
        $vals2[$c] = $ws.Cells.Item($row2, $c).Value2
    }
    for ($c = $colStart; $c -le $colEnd; $c++) {
        $ws.Cells.Item($row1, $c).Value2 = $vals2[$c]
        $ws.Cells.Item($row2, $c).Value2 = $vals1[$c]
    }
}

# Columns F..V are column indexes 6..22
Swap-RowRange 19 20 6 22
Swap-RowRange 84 85 6 22

# -----------------------------------------------------------------
# 2) Three new match rows were appended at the bottom of the sheet
#    (rows 108-110), pushing dimension from A1:V107 to A1:V110.
# -----------------------------------------------------------------
$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

$newRows = @(
    @{ A=107; B="serbia"; C="super-liga"; D="2023-2024"; E=45242.6875;        F="Javor";    G=2; H="Sp. Subotica"; I=0; J=2.09; K="09/11/2023 09:13"; L=2.16; M="11/11/2023 21:33"; N=3.11; O="09/11/2023 09:13"; P=3.09; Q="12/11/2023 15:30"; R=3.2;  S="09/11/2023 09:13"; T=2.67;                U="11/11/2023 20:34"; V="https://www.betexplorer.com/football/serbia/super-liga/javor-spartak-subotica/l0cfK6lm/" },
    @{ A=108; B="serbia"; C="super-liga"; D="2023-2024"; E=45242.79166666666; F="Partizan"; G=2; H="Cukaricki";    I=1; J=1.43; K="09/11/2023 09:13"; L=1.4;  M="12/11/2023 18:54"; N=4.09; O="09/11/2023 09:13"; P=4.25; Q="12/11/2023 18:54"; R=5.75; S="09/11/2023 09:13"; T=8.140000000000001;   U="12/11/2023 18:54"; V="https://www.betexplorer.com/football/serbia/super-liga/partizan-cukaricki/jPHJQ8BJ/" },
    @{ A=109; B="serbia"; C="super-liga"; D="2023-2024"; E=45242.8125;        F="TSC";      G=1; H="Vojvodina";    I=2; J=1.91; K="09/11/2023 09:13"; L=2.09; M="12/11/2023 19:25"; N=3.28; O="09/11/2023 09:13"; P=3.33; Q="12/11/2023 19:23"; R=3.48; S="09/11/2023 09:13"; T=3.41;                U="12/11/2023 19:25"; V="https://www.betexplorer.com/football/serbia/super-liga/tsc-backa-topola-vojvodina/29Oy2PYI/" }
)

$startRow = 108
$r = $startRow
foreach ($row in $newRows) {
    for ($ci = 0; $ci -lt $cols.Length; $ci++) {
        $col = $cols[$ci]
        $ws.Cells.Item($r, $ci + 1).Value2 = $row[$col]
    }
    $r = $r + 1
}

# Copy the formatting (bold/border style for column A, datetime style
# for column E) from the last pre-existing data row (107) down onto
# the three freshly appended rows so the new cells match the rest of
# the table's look (s="1" on A, s="2" on E).
$ws.Range("A107").Copy() | Out-Null
$ws.Range("A108:A110").PasteSpecial(-4122) | Out-Null
$ws.Range("E107").Copy() | Out-Null
$ws.Range("E108:E110").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
